# Applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed cell we force a Text number format before assigning
# the value (many of the Price column values look like plain numbers,
# e.g. "0.9991", and would otherwise be auto-converted by Excel into a
# numeric value, losing the original formatted-text representation).
# The style is reset back to "Normal" afterwards so we do not leave any
# unintended cell-formatting change behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

# Row 2
Set-CellText 'D2' '30.282.52'
Set-CellText 'E2' '  -2.09%  '

# Row 3
Set-CellText 'D3' '1.881.83'
Set-CellText 'E3' '  -2.37%  '

# Row 4
Set-CellText 'D4' '0.9991'
Set-CellText 'E4' '  -0.08%  '

# Row 5
Set-CellText 'D5' '236.56'
Set-CellText 'E5' '  -1.77%  '

# Row 6
Set-CellText 'D6' '0.9991'
Set-CellText 'E6' '  -0.07%  '

# Row 7
Set-CellText 'D7' '0.4840'
Set-CellText 'E7' '  -1.55%  '

# Row 8
Set-CellText 'D8' '0.2881'
Set-CellText 'E8' '  -3.31%  '

# Row 9
Set-CellText 'D9' '0.06602'
Set-CellText 'E9' '  -2.77%  '

# Row 10
Set-CellText 'D10' '1.881.87'
Set-CellText 'E10' '  -2.34%  '

# Row 11
Set-CellText 'D11' '16.88'
Set-CellText 'E11' '  -1.88%  '

# Row 12
Set-CellText 'D12' '0.07329'
Set-CellText 'E12' '  +0.17%  '

# Row 13
Set-CellText 'D13' '5.143'
Set-CellText 'E13' '  -0.88%  '

# Row 14
Set-CellText 'D14' '87.44'
Set-CellText 'E14' '  -2.99%  '

# Row 15
Set-CellText 'D15' '0.6583'
Set-CellText 'E15' '  -2.99%  '

# Row 16
Set-CellText 'D16' '30.258.70'
Set-CellText 'E16' '  -2.08%  '

# Row 17
Set-CellText 'D17' '13.36'
Set-CellText 'E17' '  -1.70%  '

# Row 18
Set-CellText 'B18' 'ShibaInu'
Set-CellText 'C18' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText 'D18' '0.000007758'
Set-CellText 'E18' '  -3.44%  '

# Row 19
Set-CellText 'B19' 'Dai'
Set-CellText 'C19' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText 'D19' '0.9990'
Set-CellText 'E19' '  -0.08%  '

# Row 20
Set-CellText 'D20' '5.421'
Set-CellText 'E20' '  +4.20%  '

# Row 21
Set-CellText 'D21' '2.122.87'
Set-CellText 'E21' '  -1.44%  '

# Row 22
Set-CellText 'D22' '0.9987'
Set-CellText 'E22' '  -0.10%  '

# Row 23
Set-CellText 'D23' '195.98'
Set-CellText 'E23' '  -5.38%  '

# Row 24
Set-CellText 'D24' '6.144'
Set-CellText 'E24' '  -2.62%  '

# Row 25
Set-CellText 'D25' '9.275'
Set-CellText 'E25' '  -4.58%  '

# Row 26
Set-CellText 'D26' '164.07'
Set-CellText 'E26' '  +3.01%  '

# Row 27
Set-CellText 'D27' '18.14'
Set-CellText 'E27' '  -5.07%  '

# Row 28
Set-CellText 'D28' '1.925'
Set-CellText 'E28' '  -3.97%  '

# Row 29
Set-CellText 'D29' '1.433'
Set-CellText 'E29' '  +0.42%  '

# Row 30
Set-CellText 'D30' '4.292'
Set-CellText 'E30' '  -1.89%  '

# Row 31
Set-CellText 'D31' '0.09153'
Set-CellText 'E31' '  -0.53%  '

# Row 32
Set-CellText 'D32' '4.023'
Set-CellText 'E32' '  -2.09%  '

# Row 33
Set-CellText 'E33' '  -2.85%  '

# Row 34
Set-CellText 'D34' '0.7204'
Set-CellText 'E34' '  -5.25%  '

# Row 35
Set-CellText 'D35' '1.124'
Set-CellText 'E35' '  -0.84%  '

# Row 36
Set-CellText 'D36' '2.696'
Set-CellText 'E36' '  -0.97%  '

# Row 37
Set-CellText 'D37' '0.01780'
Set-CellText 'E37' '  -4.80%  '

# Row 38
Set-CellText 'D38' '2.641'
Set-CellText 'E38' '  -3.72%  '

# Row 39
Set-CellText 'D39' '0.9189'
Set-CellText 'E39' '  -1.17%  '

# Row 40
Set-CellText 'D40' '2.047'
Set-CellText 'E40' '  -2.74%  '

# Row 41
Set-CellText 'D41' '105.82'
Set-CellText 'E41' '  -2.69%  '

# Row 42
Set-CellText 'D42' '0.4301'
Set-CellText 'E42' '  -5.44%  '

# Row 43
Set-CellText 'D43' '5.798'
Set-CellText 'E43' '  -2.10%  '

# Row 44
Set-CellText 'D44' '0.9997'
Set-CellText 'E44' '  -1.30%  '

# Row 45
Set-CellText 'D45' '7.449'
Set-CellText 'E45' '  -4.05%  '

# Row 46
Set-CellText 'D46' '0.1310'
Set-CellText 'E46' '  -6.64%  '

# Row 47
Set-CellText 'D47' '65.13'
Set-CellText 'E47' '  -7.89%  '

# Row 48
Set-CellText 'D48' '1.521'
Set-CellText 'E48' '  +5.15%  '

# Row 49
Set-CellText 'D49' '8.904'
Set-CellText 'E49' '  -2.24%  '

# Row 50
Set-CellText 'D50' '33.85'
Set-CellText 'E50' '  -5.13%  '

# Row 51
Set-CellText 'E51' '  -3.65%  '
